$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values for rows 2-5 (A/B columns)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 297

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 167

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 144

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 106

# Delete row 6 entirely (shifts cells up, removes the row)
$ws.Range("A6:B6").Delete() | Out-Null
